$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3940.5
$ws.Range("I62").Value = 3940.5
$ws.Range("K62").Value = 3940.5
$ws.Range("M62").Value = -3316.5
$ws.Range("H65").Value = 3940.5
$ws.Range("I65").Value = 3940.5
$ws.Range("K65").Value = 19702.5
$ws.Range("M65").Value = -16582.5
$ws.Range("H94").Value = 104.5
$ws.Range("J94").Value = 104
$ws.Range("L94").Value = 104
$ws.Range("N94").Value = -1006
$ws.Range("H137").Value = 1662.0333
$ws.Range("I137").Value = 1494.6538
$ws.Range("K137").Value = 4483.9614
$ws.Range("M137").Value = -1933.9614
$ws.Range("H138").Value = 1991.5454
$ws.Range("I138").Value = 1764.72
$ws.Range("J138").Value = 2180.5667
$ws.Range("K138").Value = 5294.16
$ws.Range("L138").Value = 6541.7001
$ws.Range("M138").Value = -154.1599999999999
$ws.Range("N138").Value = -16821.7001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1077.92
$ws.Range("I32").Value = 1005.6458
$ws.Range("J32").Value = 2812.5
$ws.Range("K32").Value = 1005.6458
$ws.Range("L32").Value = 2812.5
$ws.Range("M32").Value = -718.6458
$ws.Range("N32").Value = -3386.5
$ws.Range("H74").Value = 854.14703
$ws.Range("I74").Value = 763.0333000000001
$ws.Range("J74").Value = 1537.5
$ws.Range("K74").Value = 763.0333000000001
$ws.Range("L74").Value = 1537.5
$ws.Range("M74").Value = 110.9666999999999
$ws.Range("N74").Value = -3285.5
$ws.Range("H77").Value = 854.14703
$ws.Range("I77").Value = 763.0333000000001
$ws.Range("J77").Value = 1537.5
$ws.Range("K77").Value = 3815.1665
$ws.Range("L77").Value = 7687.5
$ws.Range("M77").Value = 552.8334999999997
$ws.Range("N77").Value = -16423.5
$ws.Range("H132").Value = 3427.2104
$ws.Range("I132").Value = 3664.7673
$ws.Range("J132").Value = 2697.5715
$ws.Range("K132").Value = 10994.3019
$ws.Range("L132").Value = 8092.7145
$ws.Range("M132").Value = -8464.3019
$ws.Range("N132").Value = -13152.7145
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3030.0667
$ws.Range("I20").Value = 3346.2354
$ws.Range("J20").Value = 2616.6155
$ws.Range("K20").Value = 3346.2354
$ws.Range("L20").Value = 2616.6155
$ws.Range("M20").Value = -3099.2354
$ws.Range("N20").Value = -3110.6155
$ws.Range("H86").Value = 3697.7144
$ws.Range("I86").Value = 2684.2
$ws.Range("J86").Value = 6231.5
$ws.Range("K86").Value = 2684.2
$ws.Range("L86").Value = 6231.5
$ws.Range("M86").Value = -1561.2
$ws.Range("N86").Value = -8477.5
$ws.Range("H89").Value = 3697.7144
$ws.Range("I89").Value = 2684.2
$ws.Range("J89").Value = 6231.5
$ws.Range("K89").Value = 13421
$ws.Range("L89").Value = 31157.5
$ws.Range("M89").Value = -7805
$ws.Range("N89").Value = -42389.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8130928.5
$ws.Range("I58").Value = 899
$ws.Range("J58").Value = 33334020
$ws.Range("K58").Value = 899
$ws.Range("L58").Value = 33334020
$ws.Range("M58").Value = -696
$ws.Range("N58").Value = -33334426
$ws.Range("H105").Value = 650.5714
$ws.Range("I105").Value = 472.45456
$ws.Range("J105").Value = 1303.6666
$ws.Range("K105").Value = 472.45456
$ws.Range("L105").Value = 1303.6666
$ws.Range("M105").Value = 1274.54544
$ws.Range("N105").Value = -4797.6666
$ws.Range("H136").Value = 8130928.5
$ws.Range("I136").Value = 899
$ws.Range("J136").Value = 33334020
$ws.Range("K136").Value = 2697
$ws.Range("L136").Value = 100002060
$ws.Range("M136").Value = -147
$ws.Range("N136").Value = -100007160
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 41.8
$ws.Range("J10").Value = 50
$ws.Range("L10").Value = 150
$ws.Range("N10").Value = -428
$ws.Range("H13").Value = 789.1111
$ws.Range("I13").Value = 525.5
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 1576.5
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -1408.5
$ws.Range("N13").Value = -3336
$ws.Range("H15").Value = 833.3333
$ws.Range("J15").Value = 1000
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3280
$ws.Range("H16").Value = 500
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -1846
$ws.Range("H17").Value = 436.6154
$ws.Range("I17").Value = 233.33333
$ws.Range("J17").Value = 497.6
$ws.Range("K17").Value = 699.99999
$ws.Range("L17").Value = 1492.8
$ws.Range("M17").Value = -530.99999
$ws.Range("N17").Value = -1830.8
$ws.Range("H26").Value = 33333550
$ws.Range("I26").Value = 900
$ws.Range("J26").Value = 40000080
$ws.Range("K26").Value = 2700
$ws.Range("L26").Value = 120000240
$ws.Range("M26").Value = -2412
$ws.Range("N26").Value = -120000816
$ws.Range("H29").Value = 64.666664
$ws.Range("J29").Value = 64.666664
$ws.Range("L29").Value = 193.999992
$ws.Range("N29").Value = -747.999992
$ws.Range("H34").Value = 1750.625
$ws.Range("J34").Value = 1715
$ws.Range("L34").Value = 5145
$ws.Range("N34").Value = -5313
$ws.Range("H36").Value = 1058.6428
$ws.Range("J36").Value = 999
$ws.Range("L36").Value = 2997
$ws.Range("N36").Value = -3335
$ws.Range("H44").Value = 348.625
$ws.Range("I44").Value = 384.2857
$ws.Range("J44").Value = 99
$ws.Range("K44").Value = 1152.8571
$ws.Range("L44").Value = 297
$ws.Range("M44").Value = -754.8571000000002
$ws.Range("N44").Value = -1093
$ws.Range("H131").Value = 296466.53
$ws.Range("I131").Value = 5486.364
$ws.Range("J131").Value = 381820.72
$ws.Range("K131").Value = 16459.092
$ws.Range("L131").Value = 1145462.16
$ws.Range("M131").Value = -11419.092
$ws.Range("N131").Value = -1155542.16
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 829.6667
$ws.Range("I97").Value = 724.6429000000001
$ws.Range("J97").Value = 2300
$ws.Range("K97").Value = 724.6429000000001
$ws.Range("L97").Value = 2300
$ws.Range("M97").Value = -228.6429000000001
$ws.Range("N97").Value = -3292
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1500
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1772
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 666.6667
$ws.Range("J107").Value = 933.3333
$ws.Range("K107").Value = 2000.0001
$ws.Range("L107").Value = 2799.9999
$ws.Range("M107").Value = -80.00009999999997
$ws.Range("N107").Value = -6639.9999
